$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Data": append two new weekly observations (rows 98 and 99)
# ---------------------------------------------------------------------------
$wsData = $wb.Worksheets.Item("Data")

$wsData.Range("A98").Value2 = 45147
$wsData.Range("B98").Value2 = 3222.856
$wsData.Range("A99").Value2 = 45154
$wsData.Range("B99").Value2 = 3245.971

# Match the date formatting/style already used by the preceding data rows.
$wsData.Range("A97").Copy()
$wsData.Range("A98:A99").PasteSpecial(-4122)
$wsData.Application.CutCopyMode = $false

# ---------------------------------------------------------------------------
# Sheet "SeriesInfo": refresh the metadata pulled along with the new data
# ---------------------------------------------------------------------------
$wsInfo = $wb.Worksheets.Item("SeriesInfo")

# realtime_start / realtime_end - keep these as literal text, not dates
$wsInfo.Range("B3").NumberFormat = "@"
$wsInfo.Range("B3").Value = "2023-08-22"
$wsInfo.Range("B3").ClearFormats()

$wsInfo.Range("B4").NumberFormat = "@"
$wsInfo.Range("B4").Value = "2023-08-22"
$wsInfo.Range("B4").ClearFormats()

# observation_end
$wsInfo.Range("B7").NumberFormat = "@"
$wsInfo.Range("B7").Value = "2023-08-16"
$wsInfo.Range("B7").ClearFormats()

# last_updated
$wsInfo.Range("B14").Value = "2023-08-17 15:35:51-05"

# popularity
$wsInfo.Range("B15").Value = 75
